# Bug Metrics v1.xlsx - "Debugged bug 1 and 3"
# Updates to the "Iteration 2" sheet: bug #1 (Edit Practical group account) and
# bug #3 (Delete Case Scenario (Admin)) are now marked Solved, with the fix
# description, the person who solved it, and the date solved filled in. Bug #1's
# "Date Found" value is also corrected from the mistyped text "017-10-2014" to a
# proper date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration 2")

# --- Row 8 : Bug 1 - Edit Practical group account ---
# Fix the mistyped "Date Found" (was text "017-10-2014") -> 17 Oct 2014
$ws.Range("F8").Value = 41929
# Status: Unsolved -> Solved
$ws.Range("I8").Value = "Solved"
# Solved by
$ws.Range("K8").Value = "Shi Qi"
# Action Taken by Developers
$ws.Range("J8").Value = "Changed the textbox in editPracticalGroupAccount to dropdown to eliminate error of typing the lecturer ID that does not exist"
# Date Solved -> 19 Oct 2014
$ws.Range("L8").Value = 41931

# --- Row 10 : Bug 3 - Delete Case Scenario (Admin) ---
# Status: Unsolved -> Solved
$ws.Range("I10").Value = "Solved"
# Solved by
$ws.Range("K10").Value = "Shi Qi"
# Action Taken by Developers
$ws.Range("J10").Value = "Displayed success message. Previously, the success message was set, but was not displayed on viewScenarioAdmin"
# Date Solved -> 19 Oct 2014
$ws.Range("L10").Value = 41931

# Update the row heights to fit the newly-filled wrapped text
$ws.Rows.Item(8).RowHeight = 68.25
$ws.Rows.Item(10).RowHeight = 54.75

# Restore the last active selection on this sheet
$ws.Activate()
$ws.Range("I8").Select()
